$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data: row 15 PRIORIDADE value 10 -> 5 ---
$ws.Range("C15").Value = 5

# --- Add new requisito row (row 16) ---
$ws.Range("B16").Value = "FUNCIONALIDADE DE MONITORAMENTO DA PRODUÇÃO"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7

# Copy the number formatting (centered alignment) used by the rest of the
# table's PRIORIDADE/SPRINT columns down onto the new row.
$ws.Range("C15:D15").Copy()
$ws.Range("C16:D16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Grow the "Tabela2" table so the new row becomes part of it ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B3:E16"))

# --- Extend the CHECK column conditional formatting down to row 16 ---
$ws.Range("E4:E15").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E4:E16"))

# --- Extend the CHECK column data validation (dropdown list) down to row 16 ---
$ws.Range("E4:E15").Validation.Delete()
$ws.Range("E4:E16").Validation.Add(3, 1, 1, '$G$4:$G$6')

# --- Update the active selection to match the saved workbook state ---
[void]$ws.Range("H23").Select()

Write-Host "Edit applied"
